# Add LiPo battery pack to BOM (new row 49 on Sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start the new row as a copy of the last existing row so it inherits the
# same column formatting (fills/borders/number format) used throughout the
# table, then set its height to match the other data rows.
$ws.Range("A48:K48").Copy()
$ws.Range("A49:K49").PasteSpecial(-4122)
$ws.Rows(49).RowHeight = 17

# Fill in the new part's data.
$ws.Range("H49").Value = "354"
$ws.Range("G49").Value = "Adafruit"
$ws.Range("B49").Value = "3.7V 4400mAh"
$ws.Range("F49").Value = "Lithium-Ion Battery Pack with protection circuit"
$ws.Range("J49").Value = "485-354"
$ws.Range("I49").Value = "1528-1834-ND"

$ws.Range("A49").Value = 1
$ws.Range("C49").Value = ""
$ws.Range("D49").Value = ""
$ws.Range("E49").Value = ""
$ws.Range("K49").Value = ""

# MFN / OC-DIGIKEY columns are italicised for this part, matching the style
# already used elsewhere in the sheet for similar notes.
$ws.Range("H49:I49").Font.Italic = $true

# Extend the sheet's remembered sort range (Data > Sort) to include the new
# row, same as Excel does when a sorted table grows.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("E2:E49"))
$ws.Sort.SetRange($ws.Range("A2:K49"))
$ws.Sort.Apply()
